$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Substances")

$ws.Range("A2").Value = "CPF"
$ws.Range("B2").Value = "Chlorpyrifos"
$ws.Range("C2").Value = "Organophosphate insecticide (synthetic test data)"
$ws.Range("D2").Value = 0.005
$ws.Range("E2").Value = 0.001
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "mgPerKg"
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 350.6
$ws.Range("J2").Value = $true
